$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph "aggiungiCompito(mansione: MansioneDiCucina)"
#    -> "aggiungiCompito(ricetta: Ricetta)"
#    The underlined "mansione" becomes "ricetta"; ": MansioneDiCucina"
#    is split so that ": " stays one run and "MansioneDiCucina" becomes
#    "Ricetta" in its own run (distinct from the following ")" run).
# ------------------------------------------------------------------

$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*aggiungiCompito(mansione*MansioneDiCucina*") {
        $target1 = $cand
        break
    }
}

if ($target1 -ne $null) {
    # Replace the underlined word "mansione" -> "ricetta" (format-scoped via Find in this paragraph)
    $rFind = $target1.Range
    $rFind.Find.Execute("mansione", $true, $false, $false, $false, $false, $true, 1, $false, "ricetta", 2) | Out-Null

    # Locate "MansioneDiCucina" text dynamically and replace its characters with "Ricetta"
    $full = $target1.Range
    $paraStart = $full.Start
    $relIdx = $full.Text.IndexOf("MansioneDiCucina")
    if ($relIdx -ge 0) {
        $wordStart = $paraStart + $relIdx
        $wordEnd = $wordStart + ("MansioneDiCucina".Length)
        $wordRange = $d.Range($wordStart, $wordEnd)
        $wordRange.Text = "Ricetta"
        $newWordEnd = $wordStart + ("Ricetta".Length)

        # Force run boundaries (so the new text doesn't get absorbed into the
        # neighbouring identically-formatted runs) by briefly planting
        # bookmarks at the desired split points, then removing them again.
        $b1 = $d.Range($wordStart, $wordStart)
        $d.Bookmarks.Add("zzSplitBefore", $b1) | Out-Null
        $b2 = $d.Range($newWordEnd, $newWordEnd)
        $d.Bookmarks.Add("zzSplitAfter", $b2) | Out-Null

        $d.Bookmarks.Item("zzSplitBefore").Delete()
        $d.Bookmarks.Item("zzSplitAfter").Delete()
    }
}

# ------------------------------------------------------------------
# 2) Paragraph "c consiste in mansione" -> "c consiste in ricetta"
# ------------------------------------------------------------------

$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*consiste in mansione*") {
        $target2 = $cand
        break
    }
}

if ($target2 -ne $null) {
    $rFind2 = $target2.Range
    $rFind2.Find.Execute("mansione", $true, $false, $false, $false, $false, $true, 1, $false, "ricetta", 2) | Out-Null
}
